$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 - this shifts the existing rows 6..50
# (and all their data/formatting) down to rows 7..51, exactly reproducing
# the "every row moves down by one, oldest entry lands at the new last
# row 51" pattern seen in the diff.
$ws.Rows.Item(6).Insert()

# Populate the newly-inserted row 6 with the new weekly price entry.
# Columns A,B,C,E,F,G,H,I,J,N,O,Q,R keep the same boilerplate values the
# series always uses; D (fecha), K/L/M (precios) and P (precio $/Kg) are
# the genuinely new data points.
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 45050
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112026
$ws.Range("G6").Value = "Haba"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 22000
$ws.Range("M6").Value = 21000
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 840
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
